# Updates cryptocurrency price/volume data in the "cryptos" worksheet.
# Mirrors a refreshed scrape: most rows keep the same coin/link but get
# new Price (D) and Volume(1h) (E) figures; rows 38 and 40 swap which
# coin (Mantle / Filecoin) occupies that slot, with updated price/volume.
#
# Column D occasionally holds values that look numeric (e.g. "48.82",
# "141.20"); a leading apostrophe forces Excel to store/keep them as
# text (matching the workbook's existing inlineStr/text convention and
# preserving exact formatting like trailing zeros) instead of silently
# converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.122.91"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "3.739.99"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'601.59"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'167.52"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "3.739.16"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("D11").Value = "'6.38"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "4.367.14"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "3.736.90"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "69.069.30"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "'7.37"
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("D19").Value = "'17.39"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = "  +11.19%  "
$ws.Range("D22").Value = "'492.16"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'0.729"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D25").Value = "'84.74"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "'12.28"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'2.98"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("D32").Value = "'2.47"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "'31.57"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "3.885.06"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "3.673.27"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'5.95"
$ws.Range("E38").Value = "  +2.53%  "
$ws.Range("E39").Value = "  +5.65%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'1.01"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "'0.327"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'3.04"
$ws.Range("E42").Value = "  +6.21%  "
$ws.Range("D43").Value = "'48.82"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").Value = "'423.85"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D48").Value = "'40.05"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "'141.20"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").Value = "2.780.85"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("E51").Value = "  +0.19%  "
